$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value  = "Shivamogga (Shimoga)"
$ws.Range("G8").Value  = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G11").Value = "Vijayapura (Bijapur)"
$ws.Range("G12").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G19").Value = "Bidar"
$ws.Range("G20").Value = "Ballari (Bellary)"
$ws.Range("G25").Value = "Kalaburagi (Gulbarga)"
